# Apply the final-revision edits to the "harp expander v1" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("harp expander v1")

# E6 previously held a plain numeric value (1301.9306999999999); replace it
# with the text MPN "B3F-4000" (new shared string), matching the pattern
# used by the other rows in the MPN column (E7:E10).
$ws.Range("E6").Value = "B3F-4000"

# Update the active selection to H11 (matches the saved sheetView selection).
$ws.Range("H11").Select()
